# Refresh the cryptos snapshot table (Coin / Link / Price / Volume(1h)) on Sheet1.
# D-column prices that would otherwise auto-parse as numbers are written with a
# leading apostrophe so Excel keeps them as plain text (matches the source data,
# which stores prices like "355.61" or "0.578" as strings, not numeric values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    # Force literal text even when $text looks like a number (e.g. "19.60"),
    # the same way a user typing ' before a value keeps it as text in Excel.
    $ws.Range($range).Value = '''' + $text
}

$ws.Range("D2").Value = '51.250.68'
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").Value = '2.748.70'
$ws.Range("E3").Value = '  -3.08%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextCell "D5" '355.61'
$ws.Range("E5").Value = '  -1.69%  '
Set-TextCell "D6" '107.85'
$ws.Range("E6").Value = '  -3.63%  '
$ws.Range("E7").Value = '  -2.94%  '
$ws.Range("E8").Value = '  +0.08%  '
Set-TextCell "D9" '0.578'
$ws.Range("E9").Value = '  -4.07%  '
Set-TextCell "D10" '39.10'
$ws.Range("E10").Value = '  -4.37%  '
$ws.Range("E11").Value = '  +2.96%  '
Set-TextCell "D12" '0.0834'
$ws.Range("E12").Value = '  -3.64%  '
Set-TextCell "D13" '19.60'
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("E14").Value = '  -4.07%  '
$ws.Range("D15").Value = '3.184.80'
$ws.Range("E15").Value = '  -3.13%  '
$ws.Range("D16").Value = '2.747.17'
$ws.Range("E16").Value = '  -3.45%  '
Set-TextCell "D17" '0.920'
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").Value = '51.175.31'
$ws.Range("E18").Value = '  -1.84%  '
Set-TextCell "D19" '7.54'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("E20").Value = '  -4.50%  '
Set-TextCell "D21" '12.97'
$ws.Range("E21").Value = '  -3.13%  '
$ws.Range("D22").Value = '0.0₃0957'
$ws.Range("E22").Value = '  -4.21%  '
Set-TextCell "D23" '69.34'
$ws.Range("E23").Value = '  -1.45%  '
Set-TextCell "D24" '264.59'
$ws.Range("E24").Value = '  -2.79%  '
Set-TextCell "D25" '2.73'
$ws.Range("E25").Value = '  -2.69%  '
$ws.Range("E26").Value = '  +0.03%  '
Set-TextCell "D27" '25.97'
$ws.Range("E27").Value = '  -3.47%  '
Set-TextCell "D28" '0.161'
$ws.Range("E28").Value = '  +12.95%  '
$ws.Range("E29").Value = '  +0.51%  '
Set-TextCell "D30" '10.05'
$ws.Range("E30").Value = '  -2.57%  '
Set-TextCell "D31" '34.77'
$ws.Range("E31").Value = '  -0.81%  '
Set-TextCell "D32" '6.08'
$ws.Range("E32").Value = '  +3.33%  '
Set-TextCell "D33" '51.30'
$ws.Range("E33").Value = '  -2.11%  '
Set-TextCell "D34" '0.0439'
$ws.Range("E34").Value = '  -8.19%  '
$ws.Range("E35").Value = '  -2.73%  '
Set-TextCell "D36" '5.12'
$ws.Range("E36").Value = '  -8.17%  '
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("E39").Value = '  -4.71%  '
Set-TextCell "D40" '1.94'
$ws.Range("E40").Value = '  -5.02%  '
$ws.Range("E41").Value = '  -3.18%  '
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("E43").Value = '  -2.88%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell "D44" '119.22'
$ws.Range("E44").Value = '  -4.81%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell "D45" '21.74'
$ws.Range("E45").Value = '  -3.60%  '
$ws.Range("D46").Value = '2.080.73'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell "D47" '3.21'
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell "D48" '2.29'
$ws.Range("E48").Value = '  -0.99%  '
Set-TextCell "D49" '0.923'
$ws.Range("E49").Value = '  -4.63%  '
$ws.Range("E50").Value = '  -6.65%  '
Set-TextCell "D51" '1.28'
$ws.Range("E51").Value = '  +4.18%  '
